$p = $ppt.ActivePresentation

# ----------------------------------------------------------------------
# Helper: update the cached text of a placeholder shape (if any) found
# in the given Shapes collection whose current text equals $oldText,
# replacing it with $newText.
# ----------------------------------------------------------------------
function Update-ShapeText($shapes, $oldText, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldText) {
                $tr.Text = $newText
            }
        }
    }
}

$oldDate = "16/3/2020"
$newDate = "19/3/2020"

# Slide master date placeholder
Update-ShapeText $p.SlideMaster.Shapes $oldDate $newDate

# Every slide layout's date placeholder
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-ShapeText $layout.Shapes $oldDate $newDate
}

# ----------------------------------------------------------------------
# Slide text updates
# ----------------------------------------------------------------------

# Slide 1 title: "Your slide" -> "Influence Flower"
Update-ShapeText $p.Slides.Item(1).Shapes "Your slide" "Influence Flower"

# Slide 4 title: "Case" -> "Findings"
Update-ShapeText $p.Slides.Item(4).Shapes "Case" "Findings"
